$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the whole B2:D9 block to 0 first
$ws.Range("B2:D9").Value = 0

# Then apply the specific non-zero overrides from the diff
$ws.Range("C5").Value = -0.6672883717503439
$ws.Range("C9").Value = -0.7370391310428605
